$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 239.09091
$ws.Cells.Item(8, 9).Value = 163
$ws.Cells.Item(8, 10).Value = 1000
$ws.Cells.Item(8, 11).Value = 489
$ws.Cells.Item(8, 12).Value = 3000
$ws.Cells.Item(8, 13).Value = -350
$ws.Cells.Item(8, 14).Value = -3278
$ws.Cells.Item(137, 8).Value = 813.7586
$ws.Cells.Item(137, 9).Value = 768.1053000000001
$ws.Cells.Item(137, 10).Value = 900.5
$ws.Cells.Item(137, 11).Value = 2304.3159
$ws.Cells.Item(137, 12).Value = 2701.5
$ws.Cells.Item(137, 13).Value = 245.6840999999999
$ws.Cells.Item(137, 14).Value = -7801.5
$ws.Cells.Item(141, 8).Value = 7562.727
$ws.Cells.Item(141, 9).Value = 11917
$ws.Cells.Item(141, 10).Value = 3934.1667
$ws.Cells.Item(141, 11).Value = 35751
$ws.Cells.Item(141, 12).Value = 11802.5001
$ws.Cells.Item(141, 13).Value = -30571
$ws.Cells.Item(141, 14).Value = -22162.5001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1172.7333
$ws.Cells.Item(61, 9).Value = 1035.9524
$ws.Cells.Item(61, 10).Value = 1491.8889
$ws.Cells.Item(61, 11).Value = 1035.9524
$ws.Cells.Item(61, 12).Value = 1491.8889
$ws.Cells.Item(61, 13).Value = -823.9523999999999
$ws.Cells.Item(61, 14).Value = -1915.8889
$ws.Cells.Item(74, 8).Value = 1221.4286
$ws.Cells.Item(74, 9).Value = 1110.7333
$ws.Cells.Item(74, 10).Value = 1498.1666
$ws.Cells.Item(74, 11).Value = 1110.7333
$ws.Cells.Item(74, 12).Value = 1498.1666
$ws.Cells.Item(74, 13).Value = -236.7333000000001
$ws.Cells.Item(74, 14).Value = -3246.1666
$ws.Cells.Item(77, 8).Value = 1221.4286
$ws.Cells.Item(77, 9).Value = 1110.7333
$ws.Cells.Item(77, 10).Value = 1498.1666
$ws.Cells.Item(77, 11).Value = 5553.6665
$ws.Cells.Item(77, 12).Value = 7490.833000000001
$ws.Cells.Item(77, 13).Value = -1185.6665
$ws.Cells.Item(77, 14).Value = -16226.833
$ws.Cells.Item(110, 8).Value = 1400
$ws.Cells.Item(110, 9).Value = 1400
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 11).Value = 1400
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 13).Value = 645
$ws.Cells.Item(130, 8).Value = 25499.75
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 25499.75
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 25499.75
$ws.Cells.Item(130, 14).Value = -35539.75
$ws.Cells.Item(132, 8).Value = 71502616
$ws.Cells.Item(132, 9).Value = 333334620
$ws.Cells.Item(132, 10).Value = 93892
$ws.Cells.Item(132, 11).Value = 1000003860
$ws.Cells.Item(132, 12).Value = 281676
$ws.Cells.Item(132, 13).Value = -1000001330
$ws.Cells.Item(132, 14).Value = -286736
$ws.Cells.Item(136, 8).Value = 1172.7333
$ws.Cells.Item(136, 9).Value = 1035.9524
$ws.Cells.Item(136, 10).Value = 1491.8889
$ws.Cells.Item(136, 11).Value = 3107.857199999999
$ws.Cells.Item(136, 12).Value = 4475.6667
$ws.Cells.Item(136, 13).Value = -557.8571999999995
$ws.Cells.Item(136, 14).Value = -9575.6667

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 8079.2334
$ws.Cells.Item(134, 9).Value = 3408.1667
$ws.Cells.Item(134, 10).Value = 15085.833
$ws.Cells.Item(134, 11).Value = 10224.5001
$ws.Cells.Item(134, 12).Value = 45257.499
$ws.Cells.Item(134, 13).Value = -7689.500100000001
$ws.Cells.Item(134, 14).Value = -50327.499

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6581173.5
$ws.Cells.Item(31, 9).Value = 8066303
$ws.Cells.Item(31, 10).Value = 4171.4287
$ws.Cells.Item(31, 11).Value = 8066303
$ws.Cells.Item(31, 12).Value = 4171.4287
$ws.Cells.Item(31, 13).Value = -8066008
$ws.Cells.Item(31, 14).Value = -4761.4287
$ws.Cells.Item(34, 8).Value = 6581173.5
$ws.Cells.Item(34, 9).Value = 8066303
$ws.Cells.Item(34, 10).Value = 4171.4287
$ws.Cells.Item(34, 11).Value = 8066303
$ws.Cells.Item(34, 12).Value = 4171.4287
$ws.Cells.Item(34, 13).Value = -8066101
$ws.Cells.Item(34, 14).Value = -4575.4287
$ws.Cells.Item(58, 8).Value = 1129.1578
$ws.Cells.Item(58, 9).Value = 1070.3636
$ws.Cells.Item(58, 10).Value = 1210
$ws.Cells.Item(58, 11).Value = 1070.3636
$ws.Cells.Item(58, 12).Value = 1210
$ws.Cells.Item(58, 13).Value = -867.3635999999999
$ws.Cells.Item(58, 14).Value = -1616
$ws.Cells.Item(132, 8).Value = 39629.184
$ws.Cells.Item(132, 9).Value = 1674.5
$ws.Cells.Item(132, 10).Value = 148071.14
$ws.Cells.Item(132, 11).Value = 5023.5
$ws.Cells.Item(132, 12).Value = 444213.42
$ws.Cells.Item(132, 13).Value = -2493.5
$ws.Cells.Item(132, 14).Value = -449273.42
$ws.Cells.Item(134, 8).Value = 2122.6956
$ws.Cells.Item(134, 9).Value = 1399.6666
$ws.Cells.Item(134, 10).Value = 4725.6
$ws.Cells.Item(134, 11).Value = 4198.9998
$ws.Cells.Item(134, 12).Value = 14176.8
$ws.Cells.Item(134, 13).Value = -1663.9998
$ws.Cells.Item(134, 14).Value = -19246.8
$ws.Cells.Item(136, 8).Value = 1129.1578
$ws.Cells.Item(136, 9).Value = 1070.3636
$ws.Cells.Item(136, 10).Value = 1210
$ws.Cells.Item(136, 11).Value = 3211.0908
$ws.Cells.Item(136, 12).Value = 3630
$ws.Cells.Item(136, 13).Value = -661.0907999999999
$ws.Cells.Item(136, 14).Value = -8730

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 15874811
$ws.Cells.Item(129, 9).Value = 500
$ws.Cells.Item(129, 10).Value = 16668526
$ws.Cells.Item(129, 11).Value = 1500
$ws.Cells.Item(129, 12).Value = 50005578
$ws.Cells.Item(129, 13).Value = 3500
$ws.Cells.Item(129, 14).Value = -50015578

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 442265.6
$ws.Cells.Item(132, 9).Value = 125989.5
$ws.Cells.Item(132, 10).Value = 591101.4399999999
$ws.Cells.Item(132, 11).Value = 377968.5
$ws.Cells.Item(132, 12).Value = 1773304.32
$ws.Cells.Item(132, 13).Value = -375438.5
$ws.Cells.Item(132, 14).Value = -1778364.32

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 2308.7368
$ws.Cells.Item(122, 9).Value = 2250.647
$ws.Cells.Item(122, 10).Value = 2802.5
$ws.Cells.Item(122, 11).Value = 6751.941
$ws.Cells.Item(122, 12).Value = 8407.5
$ws.Cells.Item(122, 13).Value = -4301.941
$ws.Cells.Item(122, 14).Value = -13307.5
$ws.Cells.Item(128, 8).Value = 41064.832
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 41064.832
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 41064.832
$ws.Cells.Item(128, 14).Value = -51024.832
$ws.Cells.Item(132, 8).Value = 20806.115
$ws.Cells.Item(132, 9).Value = 33843.613
$ws.Cells.Item(132, 10).Value = 1560.2858
$ws.Cells.Item(132, 11).Value = 101530.839
$ws.Cells.Item(132, 12).Value = 4680.857400000001
$ws.Cells.Item(132, 13).Value = -99000.83899999999
$ws.Cells.Item(132, 14).Value = -9740.857400000001
$ws.Cells.Item(136, 8).Value = 5433.579
$ws.Cells.Item(136, 9).Value = 7187.4
$ws.Cells.Item(136, 10).Value = 3484.889
$ws.Cells.Item(136, 11).Value = 21562.2
$ws.Cells.Item(136, 12).Value = 10454.667
$ws.Cells.Item(136, 13).Value = -19012.2
$ws.Cells.Item(136, 14).Value = -15554.667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 2722.125
$ws.Cells.Item(2, 9).Value = 2267.3333
$ws.Cells.Item(2, 10).Value = 2995
$ws.Cells.Item(2, 11).Value = 2267.3333
$ws.Cells.Item(2, 12).Value = 2995
$ws.Cells.Item(2, 13).Value = -2155.3333
$ws.Cells.Item(2, 14).Value = -3219
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(132, 8).Value = 67420770
$ws.Cells.Item(132, 9).Value = 75334500
$ws.Cells.Item(132, 10).Value = 8067766
$ws.Cells.Item(132, 11).Value = 226003500
$ws.Cells.Item(132, 12).Value = 24203298
$ws.Cells.Item(132, 13).Value = -226000970
$ws.Cells.Item(132, 14).Value = -24208358
$ws.Cells.Item(136, 8).Value = 25580.902
$ws.Cells.Item(136, 9).Value = 44152.477
$ws.Cells.Item(136, 10).Value = 1850.5555
$ws.Cells.Item(136, 11).Value = 132457.431
$ws.Cells.Item(136, 12).Value = 5551.666499999999
$ws.Cells.Item(136, 13).Value = -129907.431
$ws.Cells.Item(136, 14).Value = -10651.6665

# --- Cell deletions (cells removed entirely in target) ---
$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Cells.Item(110, 14).ClearContents()   # N110 removed
$wsWVR = $wb.Worksheets.Item("WVR")
$wsWVR.Cells.Item(5, 14).ClearContents()     # N5 removed
